$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain duplicated data and both need
# their F3 (7 -> 8) and F6 (431 -> 432) cell values bumped by 1.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 8
    $ws.Range("F6").Value = 432
}
